$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (previously row 4's data)
$ws.Range("D2").Value = 44357
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 725

# Row 3 (previously row 2's data)
$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("S3").Value = 825

# Row 4 (previously row 3's data)
$ws.Range("D4").Value = 44533
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("S4").Value = 825
